$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before the existing "Comment" column (column X),
# shifting Comment (and the rest) one column to the right (to Y).
$ws.Range("X1:X5").EntireColumn.Insert("xlShiftToRight")

# Populate the newly inserted column X with the "UsedEnzyme" field.
$ws.Range("X1").Value = "UsedEnzyme"
$ws.Range("X2").Value = "# Enzyme utilisée"
$ws.Range("X3").Value = "#string"
$ws.Range("X4").Value = "# format: texte"
$ws.Range("X5").Value = "# ex:"

# Tweak the Wavelength format hint to mention the unit (nm).
$ws.Range("J4").Value = "# format: nombre entier, ne pas spécifier d'unité (nm)"
